# Applies the "added mult table modal, some lesson tweaks, prep lesson 3" edit.
# Rows 1-30 of the "en" sheet are untouched; rows 31-55 are rewritten (some
# existing rows are reshuffled/reworded, and several brand-new rows are
# appended for the new multiplication-table modal and lesson-3 prep copy).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 31-55 : new / reordered key-value (-voiceDuration) rows ---------

$ws.Range("A31").Value = "multiplicationTable"
$ws.Range("B31").Value = "Multiplication Table"

$ws.Range("A32").Value = "multTable_instruct"
$ws.Range("B32").Value = "Press this button to review the multiplication table."
$ws.Range("C32").Value = 5

$ws.Range("A33").Value = "proceed_instruct"
$ws.Range("B33").Value = "Press this button to proceed."
$ws.Range("C33").Value = 5

$ws.Range("A34").Value = "division"
$ws.Range("B34").Value = "Division"

$ws.Range("A35").Value = "commutative_title"
$ws.Range("B35").Value = "Commutative Property"

$ws.Range("A36").Value = "not_commutative"
$ws.Range("B36").Value = "Not Commutative!"

$ws.Range("A37").Value = "lesson_1_intro_1"
$ws.Range("B37").Value = "Before we proceed, let's first learn some tricks with multiplication!"

$ws.Range("A38").Value = "lesson_1_mult2_1"
$ws.Range("B38").Value = "In multiples of two, the trick is to simply double the number."

$ws.Range("A39").Value = "lesson_1_mult2_2"
$ws.Range("B39").Value = "For example: 2 x 6 can be 6 + 6, which equals to 12."

$ws.Range("A40").Value = "lesson_1_commutative_1"
$ws.Range("B40").Value = "The commutative property means that multiplying numbers in any order gives the same answer."

$ws.Range("A41").Value = "lesson_1_commutative_2"
$ws.Range("B41").Value = "For example: 2 x 3, and 3 x 2, equal 6."

$ws.Range("A42").Value = "lesson_1_commutative_3"
$ws.Range("B42").Value = "With this trick, you only have to remember half the multiplication table!"

$ws.Range("A43").Value = "lesson_1_tutorial_1"
$ws.Range("B43").Value = "Now banish these blobs by connecting them in the correct order using multiplication."

$ws.Range("A44").Value = "lesson_1_tutorial_end_1"
$ws.Range("B44").Value = "Excellent! You are now ready for the mission!"

$ws.Range("A45").Value = "lesson_2_intro_1"
$ws.Range("B45").Value = "Good work! Now it's time to step up the game with multiples of 3 and 4."

$ws.Range("A46").Value = "lesson_2_intro_2"
$ws.Range("B46").Value = "Let me show you some neat tricks."

$ws.Range("A47").Value = "lesson_2_mult3_1"
$ws.Range("B47").Value = "In multiples of three: double the number, and then add the original number."

$ws.Range("A48").Value = "lesson_2_mult3_2"
$ws.Range("B48").Value = "For example: 3 x 6, double 6 to get 12, and then add 6 to get 18."

$ws.Range("A49").Value = "lesson_2_mult4_1"
$ws.Range("B49").Value = "In multiples of four: double the number, and then double it again."

$ws.Range("A50").Value = "lesson_2_mult4_2"
$ws.Range("B50").Value = "For example: 4 x 6, double 6 to get 12, and then double 12 to get 24."

$ws.Range("A51").Value = "lesson_2_div_1"
$ws.Range("B51").Value = "When it comes to division, think of it as the opposite of multiplication."

$ws.Range("A52").Value = "lesson_2_div_2"
$ws.Range("B52").Value = "Rearranging the equation, and replacing division with multiplication can help."

$ws.Range("A53").Value = "lesson_2_div_3"
$ws.Range("B53").Value = "Unlike multiplication, division is not commutative. So the order of the numbers cannot be changed."

$ws.Range("A54").Value = "lesson_2_tutorial_1"
$ws.Range("B54").Value = "For the next mission, some blobs must be matched with division. Go ahead and try it out."

$ws.Range("A55").Value = "lesson_2_tutorial_end_1"
$ws.Range("B55").Value = "Excellent! You are now ready for the mission!"

# --- refresh the selection to match the author's final cursor ------------

$ws.Range("B55").Select()
